$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = -0.035489237677883236
$ws.Range("B1").Value = 0.035489237657426967

$ws.Range("A2").Value = 0.046977678774061955
$ws.Range("B2").Value = -0.046977678784773352

$ws.Range("A3").Value = 0.029435063627344767
$ws.Range("B3").Value = -0.029435063657541057

$ws.Range("A4").Value = 0.059235648910738584
$ws.Range("B4").Value = -0.059235648967865068

$ws.Range("A5").Value = -0.035507653392175632
$ws.Range("B5").Value = 0.035507653330246615
